$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.307.94'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '1.601.74'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.503'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0853'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = '1.826.56'
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").Value = '1.598.78'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '26.305.03'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '226.85'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.47%  '
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.60%  '
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.32%  '
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.55%  '
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0494'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("E31").Value = '  +0.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.74%  '
$ws.Range("D33").Value = '1.444.11'
$ws.Range("E33").Value = '  +8.30%  '
$ws.Range("E34").Value = '  +0.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.43'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("E36").Value = '  +1.08%  '
$ws.Range("E37").Value = '  -2.93%  '
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.824'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.83%  '
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.925'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("D44").Value = '1.738.05'
$ws.Range("E44").Value = '  +0.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.760'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("E48").Value = '  +0.70%  '
$ws.Range("D49").Value = '0.0₇0986'
$ws.Range("E49").Value = '  -3.85%  '
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0951'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.53%  '
